# Adds the "mean length of turn" (l_N, column I) data that was missing for the
# lower portion of the core table, and normalizes the formatting of a handful
# of pre-existing column-I cells whose style had drifted to a "filled"
# variant of the same font/border/alignment combo used elsewhere in the
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows whose l_N value is simply unknown -> literal "Not Available" text,
#    matching the plain (un-styled) cells already used higher up the table.
# ---------------------------------------------------------------------------
$naRows = @(71,72,73,74,75,76,77,78,93,94,104,105,106,107,108)
foreach ($r in $naRows) {
    $ws.Range("I$r").Value = "Not Available"
}

# ---------------------------------------------------------------------------
# 2) Brand-new numeric l_N values for rows 79-115 (grouped by which existing
#    cell already carries the correct font/border/alignment combination, so
#    that pasting its format reuses the matching style instead of fabricating
#    a new one).
# ---------------------------------------------------------------------------

# font19 / no border / horizontal-left  (same look as H43:H46)
$ws.Range("H43").Copy()
$ws.Range("I79").PasteSpecial(-4122)

# font18 / medium border both sides / left+wrap (same look as H80:H81)
$ws.Range("H80").Copy()
foreach ($r in @(80,81)) {
    $ws.Range("I$r").PasteSpecial(-4122)
}

# font18 / no border / left+wrap (same look as H56:H57)
$ws.Range("H56").Copy()
$group11 = @(82,83,84,85,86,87,88,89,90,91,92,95,96,97,98,99,100,101,102,103,112,113,114,115)
foreach ($r in $group11) {
    $ws.Range("I$r").PasteSpecial(-4122)
}

# font18 / left border only / left+wrap (same look as H3:H4)
$ws.Range("H3").Copy()
foreach ($r in @(109,110,111)) {
    $ws.Range("I$r").PasteSpecial(-4122)
}

$numericValues = @{
    79 = 20
    80 = 25
    81 = 30
    82 = 36
    83 = 42
    84 = 52
    85 = 61
    86 = 71.5
    87 = 96.8
    88 = 120
    89 = 140
    90 = 158
    91 = 210
    92 = 5.8
    95 = 14.6
    96 = 19.2
    97 = 22
    98 = 28
    99 = 35.6
    100 = 44
    101 = 52
    102 = 60
    103 = 73
    109 = 56
    110 = 56
    111 = 66
    112 = 62
    113 = 76
    114 = 87
    115 = 100.5
}
foreach ($r in $numericValues.Keys) {
    $ws.Range("I$r").Value = $numericValues[$r]
}

# ---------------------------------------------------------------------------
# 3) Existing column-I cells (rows 3-57) that were styled with a stray
#    "filled" variant -> repoint them at the plain variant already used by
#    the neighbouring H column, without touching their values.
# ---------------------------------------------------------------------------
$ws.Range("H3").Copy()
foreach ($r in @(3,4)) {
    $ws.Range("I$r").PasteSpecial(-4122)
}

$ws.Range("H56").Copy()
$plainGroup = @(7,11,14,18,20,21,23,27,28,31,32,33,35,36,37,38,41,56,57)
foreach ($r in $plainGroup) {
    $ws.Range("I$r").PasteSpecial(-4122)
}

$ws.Range("H43").Copy()
foreach ($r in @(43,44,45,46)) {
    $ws.Range("I$r").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Scroll position / active selection moved down to the newly-filled tail
#    of the table.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("I116").Select()
